$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "50.969.94"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.931.52"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "373.67"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.47"
$ws.Range("E6").Value = "  -4.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.535"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.581"
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.97"
$ws.Range("E10").Value = "  -3.77%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0844"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.399.70"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.93"
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.44"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "11.66"
$ws.Range("E16").Value = "  +57.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.931.01"
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.968"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "50.965.49"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("E20").Value = "  -6.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.46"
$ws.Range("E21").Value = "  -3.92%  "
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "264.45"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.46"
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.13"
$ws.Range("E25").Value = "  +10.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.00"
$ws.Range("E26").Value = "  +3.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.47"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -4.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.49"
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("E31").Value = "  -2.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.95"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.55"
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.02"
$ws.Range("E34").Value = "  -2.94%  "
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.80"
$ws.Range("E36").Value = "  -7.44%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.19"
$ws.Range("E38").Value = "  +3.86%  "
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.28"
$ws.Range("E40").Value = "  -5.68%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.46"
$ws.Range("E41").Value = "  -5.65%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.77"
$ws.Range("E42").Value = "  -4.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.99"
$ws.Range("E43").Value = "  -3.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.97"
$ws.Range("E44").Value = "  -5.09%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.04"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.274"
$ws.Range("E46").Value = "  -5.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.29"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.002.62"
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("E49").Value = "  -2.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0331"
$ws.Range("E50").Value = "  -3.95%  "
$ws.Range("E51").Value = "  +0.42%  "
